$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-26 Friday", "2024-07-27 Saturday"),
    @("47×81=3807", "73×44=3212"),
    @("64×95=6080", "13×98=1274"),
    @("69×44=3036", "15×56=840"),
    @("24×70=1680", "43×92=3956"),
    @("22×29=638", "96×96=9216"),
    @("23×47=1081", "20×46=920"),
    @("13×93=1209", "27×37=999"),
    @("70×35=2450", "12×38=456"),
    @("70×41=2870", "85×15=1275"),
    @("29×58=1682", "39×77=3003"),
    @("49×29=1421", "32×52=1664"),
    @("35×93=3255", "91×48=4368"),
    @("31×97=3007", "22×97=2134"),
    @("49×20=980", "78×62=4836"),
    @("85×41=3485", "24×82=1968"),
    @("75×29=2175", "78×50=3900"),
    @("82×79=6478", "50×31=1550"),
    @("41×53=2173", "57×98=5586"),
    @("16×56=896", "75×55=4125"),
    @("11×64=704", "19×27=513"),
    @("56×28=1568", "90×34=3060"),
    @("37×33=1221", "45×46=2070"),
    @("71×41=2911", "20×43=860"),
    @("61×57=3477", "33×97=3201"),
    @("25×48=1200", "52×24=1248")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
